$wb = $excel.ActiveWorkbook

# Sheet ALC, row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 200.5
$ws.Range("I9").Value = 200.5
$ws.Range("K9").Value = 200.5
$ws.Range("M9").Value = -31.5

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1587.1818
$ws.Range("J17").Value = 1115.9
$ws.Range("L17").Value = 3347.7
$ws.Range("N17").Value = -3683.7

# Sheet ALC, row 34
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 12643.875
$ws.Range("I34").Value = 12643.875
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 12643.875
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -12440.875
$ws.Range("N34").ClearContents()

# Sheet ALC, row 36
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 12643.875
$ws.Range("I36").Value = 12643.875
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 12643.875
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -11928.875
$ws.Range("N36").ClearContents()

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1249.6666
$ws.Range("I40").Value = 1249.5
$ws.Range("J40").Value = 1250
$ws.Range("K40").Value = 1249.5
$ws.Range("L40").Value = 1250
$ws.Range("M40").Value = -1074.5
$ws.Range("N40").Value = -1600

# Sheet ALC, row 63
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Sheet ALC, row 66
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4680
$ws.Range("I86").Value = 4351.2
$ws.Range("J86").Value = 5502
$ws.Range("K86").Value = 4351.2
$ws.Range("L86").Value = 5502
$ws.Range("M86").Value = -3228.2
$ws.Range("N86").Value = -7748

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4680
$ws.Range("I89").Value = 4351.2
$ws.Range("J89").Value = 5502
$ws.Range("K89").Value = 21756
$ws.Range("L89").Value = 27510
$ws.Range("M89").Value = -16140
$ws.Range("N89").Value = -38742

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1789.2858
$ws.Range("I112").Value = 1584.7693
$ws.Range("K112").Value = 4754.3079
$ws.Range("M112").Value = -3646.3079

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2667.8857
$ws.Range("I138").Value = 2524.3438
$ws.Range("K138").Value = 7573.0314
$ws.Range("M138").Value = -2433.0314

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3076.0833
$ws.Range("J141").Value = 3447.3333
$ws.Range("L141").Value = 10341.9999
$ws.Range("N141").Value = -20701.9999

# Sheet ARM, row 18
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 40000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 40000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 40000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -40644

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3777.9363
$ws.Range("I32").Value = 1876.8889
$ws.Range("K32").Value = 1876.8889
$ws.Range("M32").Value = -1589.8889

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7214.3706
$ws.Range("I61").Value = 5446.3
$ws.Range("J61").Value = 8254.412
$ws.Range("K61").Value = 5446.3
$ws.Range("L61").Value = 8254.412
$ws.Range("M61").Value = -5234.3
$ws.Range("N61").Value = -8678.412

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3196.8462
$ws.Range("I74").Value = 2550.6956
$ws.Range("J74").Value = 4125.6875
$ws.Range("K74").Value = 2550.6956
$ws.Range("L74").Value = 4125.6875
$ws.Range("M74").Value = -1676.6956
$ws.Range("N74").Value = -5873.6875

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3196.8462
$ws.Range("I77").Value = 2550.6956
$ws.Range("J77").Value = 4125.6875
$ws.Range("K77").Value = 12753.478
$ws.Range("L77").Value = 20628.4375
$ws.Range("M77").Value = -8385.477999999999
$ws.Range("N77").Value = -29364.4375

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2446.5
$ws.Range("I88").Value = 1944.25
$ws.Range("K88").Value = 1944.25
$ws.Range("M88").Value = -1538.25

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2446.5
$ws.Range("I91").Value = 1944.25
$ws.Range("K91").Value = 1944.25
$ws.Range("M91").Value = -540.25

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2659.2354
$ws.Range("I102").Value = 2467.1333
$ws.Range("K102").Value = 2467.1333
$ws.Range("M102").Value = -845.1333

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3366.6667
$ws.Range("I110").Value = 3244.4443
$ws.Range("J110").Value = 3733.3333
$ws.Range("K110").Value = 3244.4443
$ws.Range("L110").Value = 3733.3333
$ws.Range("M110").Value = -1199.4443
$ws.Range("N110").Value = -7823.3333

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6848
$ws.Range("I132").Value = 9524.75
$ws.Range("K132").Value = 28574.25
$ws.Range("M132").Value = -26044.25

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7214.3706
$ws.Range("I136").Value = 5446.3
$ws.Range("J136").Value = 8254.412
$ws.Range("K136").Value = 16338.9
$ws.Range("L136").Value = 24763.236
$ws.Range("M136").Value = -13788.9
$ws.Range("N136").Value = -29863.236

# Sheet BSM, row 74
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 65000
$ws.Range("J74").Value = 65000
$ws.Range("L74").Value = 65000
$ws.Range("N74").Value = -66872

# Sheet BSM, row 77
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 65000
$ws.Range("J77").Value = 65000
$ws.Range("L77").Value = 195000
$ws.Range("N77").Value = -204360

# Sheet BSM, row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 21132.375
$ws.Range("J81").Value = 22294.143
$ws.Range("L81").Value = 22294.143
$ws.Range("N81").Value = -24416.143

# Sheet BSM, row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 21132.375
$ws.Range("J84").Value = 22294.143
$ws.Range("L84").Value = 66882.429
$ws.Range("N84").Value = -77490.429

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 805.88464
$ws.Range("I94").Value = 797.875
$ws.Range("J94").Value = 902
$ws.Range("K94").Value = 797.875
$ws.Range("L94").Value = 902
$ws.Range("M94").Value = -346.875
$ws.Range("N94").Value = -1804

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3499.2104
$ws.Range("I99").Value = 2955.625
$ws.Range("J99").Value = 6398.3335
$ws.Range("K99").Value = 2955.625
$ws.Range("L99").Value = 6398.3335
$ws.Range("M99").Value = -1457.625
$ws.Range("N99").Value = -9394.333500000001

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3576.0967
$ws.Range("I105").Value = 4312.0713
$ws.Range("K105").Value = 4312.0713
$ws.Range("M105").Value = -2565.0713

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2487.818
$ws.Range("I107").Value = 2187.875
$ws.Range("J107").Value = 3287.6667
$ws.Range("K107").Value = 2187.875
$ws.Range("L107").Value = 3287.6667
$ws.Range("M107").Value = -267.875
$ws.Range("N107").Value = -7127.6667

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18721.076
$ws.Range("I134").Value = 8708.777
$ws.Range("J134").Value = 41248.75
$ws.Range("K134").Value = 26126.331
$ws.Range("L134").Value = 123746.25
$ws.Range("M134").Value = -23591.331
$ws.Range("N134").Value = -128816.25

# Sheet BSM, row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

# Sheet BSM, row 139
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280

# Sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18996.75
$ws.Range("I60").Value = 7998.5
$ws.Range("K60").Value = 7998.5
$ws.Range("M60").Value = -7487.5

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 20001598
$ws.Range("J122").Value = 25001500
$ws.Range("L122").Value = 225013500
$ws.Range("N122").Value = -225018400

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5090012
$ws.Range("I132").Value = 1264656.1
$ws.Range("K132").Value = 11381904.9
$ws.Range("M132").Value = -11379374.9

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 990.8182
$ws.Range("I139").Value = 990.8182
$ws.Range("K139").Value = 2972.4546
$ws.Range("M139").Value = 2167.5454

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 900.2093
$ws.Range("I97").Value = 752.23334
$ws.Range("J97").Value = 1241.6923
$ws.Range("K97").Value = 752.23334
$ws.Range("L97").Value = 1241.6923
$ws.Range("M97").Value = -256.23334
$ws.Range("N97").Value = -2233.6923

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 21526.916
$ws.Range("I132").Value = 24922.4
$ws.Range("J132").Value = 4549.5
$ws.Range("K132").Value = 74767.20000000001
$ws.Range("L132").Value = 13648.5
$ws.Range("M132").Value = -72237.20000000001
$ws.Range("N132").Value = -18708.5

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4995.6
$ws.Range("I93").Value = 5152.5835
$ws.Range("K93").Value = 5152.5835
$ws.Range("M93").Value = -3904.5835

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2854.3872
$ws.Range("I132").Value = 2027.2778
$ws.Range("K132").Value = 6081.8334
$ws.Range("M132").Value = -3551.8334

# Sheet LTW, row 138
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 20000
$ws.Range("J138").Value = 20000
$ws.Range("L138").Value = 20000
$ws.Range("N138").Value = -30280

# Sheet LTW, row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 59995
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 59995
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 59995
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -70355

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4830
$ws.Range("I62").Value = 2750
$ws.Range("K62").Value = 2750
$ws.Range("M62").Value = -2126

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4830
$ws.Range("I65").Value = 2750
$ws.Range("K65").Value = 13750
$ws.Range("M65").Value = -10630

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22922.521
$ws.Range("I132").Value = 14146.593
$ws.Range("K132").Value = 42439.779
$ws.Range("M132").Value = -39909.779
